# NIT-9010580831.xlsx - "Estado de Cuenta" update
# The previous periods-in-arrears (Periodo Mora) / Valor Mora rows are removed
# and replaced with the new ones: the block of 7 rows (16-22) is reversed so
# the most recent period (2003) is listed first and the oldest (1909) last,
# carrying its "Valor Mora" value (32021) along with it; the other rows keep
# the standard 33125 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora column (E16:E22) - text-formatted cells, reorder descending
$ws.Range("E16").Value = "2003"
$ws.Range("E17").Value = "2002"
$ws.Range("E18").Value = "2001"
$ws.Range("E19").Value = "1912"
$ws.Range("E20").Value = "1911"
$ws.Range("E21").Value = "1910"
$ws.Range("E22").Value = "1909"

# Valor Mora column (F16:F22) - the odd value (32021) now travels with period 2003
$ws.Range("F16").Value = 32021
$ws.Range("F17").Value = 33125
$ws.Range("F18").Value = 33125
$ws.Range("F19").Value = 33125
$ws.Range("F20").Value = 33125
$ws.Range("F21").Value = 33125
$ws.Range("F22").Value = 33125
